$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit re-derives ax..gz (cols C-H) for the existing 20 samples (rows 2-21)
# by shifting them down one row (a new leading sample lands in row 2) and appends
# 9 brand-new samples (rows 22-31, timestamps 2000..2900). Timestamp/label (A,B) for
# rows 2-21 are unchanged; new rows get their own timestamp/label too.

$data = New-Object 'object[,]' 30,8

$data[0,0] = 0.0
$data[0,1] = "struggle"
$data[0,2] = 0.2003364562988281
$data[0,3] = 0.6192827224731445
$data[0,4] = -0.0097306966781616
$data[0,5] = -0.1266735037978812
$data[0,6] = -1.426815969603404
$data[0,7] = 0.2018442184341198

$data[1,0] = 100.0
$data[1,1] = "struggle"
$data[1,2] = 0.1392936706542968
$data[1,3] = 0.3245168924331665
$data[1,4] = 0.0351060479879379
$data[1,5] = -0.1538507725511279
$data[1,6] = 0.2260856117521008
$data[1,7] = -0.09954921262604832

$data[2,0] = 200.0
$data[2,1] = "struggle"
$data[2,2] = 1.261228561401367
$data[2,3] = 1.41911768913269
$data[2,4] = -4.086977005004883
$data[2,5] = 0.4955612986671659
$data[2,6] = 3.062298653077098
$data[2,7] = -0.5691612070920504

$data[3,0] = 300.0
$data[3,1] = "struggle"
$data[3,2] = 2.267774105072021
$data[3,3] = -1.184367299079895
$data[3,4] = -3.818997621536255
$data[3,5] = 1.454008883359481
$data[3,6] = 4.839717144868811
$data[3,7] = -0.2820296159812385

$data[4,0] = 400.0
$data[4,1] = "struggle"
$data[4,2] = -1.021368980407715
$data[4,3] = 2.765882253646851
$data[4,4] = 1.220625877380371
$data[4,5] = 1.137240985218359
$data[4,6] = 3.374921506764936
$data[4,7] = 0.5055315214760452

$data[5,0] = 500.0
$data[5,1] = "struggle"
$data[5,2] = -1.566243886947632
$data[5,3] = 0.6303287744522095
$data[5,4] = 0.5539150238037109
$data[5,5] = 0.2921525373750798
$data[5,6] = 1.951630319867816
$data[5,7] = 0.7453397719227535

$data[6,0] = 600.0
$data[6,1] = "struggle"
$data[6,2] = -2.018197059631348
$data[6,3] = 1.482012748718261
$data[6,4] = 1.627924919128418
$data[6,5] = 0.1624465198541174
$data[6,6] = 1.506097605033798
$data[6,7] = 0.4439899596024537

$data[7,0] = 700.0
$data[7,1] = "struggle"
$data[7,2] = -4.581077575683594
$data[7,3] = 2.173830270767212
$data[7,4] = 9.107954978942873
$data[7,5] = 0.2823693378239257
$data[7,6] = -0.7324928641319217
$data[7,7] = -0.04955176400894934

$data[8,0] = 800.0
$data[8,1] = "struggle"
$data[8,2] = -0.2751750946044922
$data[8,3] = -0.9806771278381348
$data[8,4] = -2.06553053855896
$data[8,5] = -0.6597999164036341
$data[8,6] = -4.261745044163296
$data[8,7] = -0.7145596061434063

$data[9,0] = 900.0
$data[9,1] = "struggle"
$data[9,2] = 2.35319709777832
$data[9,3] = 1.910999298095703
$data[9,4] = -2.476747035980225
$data[9,5] = -1.235334702292262
$data[9,6] = -4.321111241165461
$data[9,7] = -0.293212206996212

$data[10,0] = 1000.0
$data[10,1] = "struggle"
$data[10,2] = -4.167366027832031
$data[10,3] = 0.0987618193030357
$data[10,4] = 3.953242778778076
$data[10,5] = 0.393702644170544
$data[10,6] = -1.912042505887086
$data[10,7] = 0.653996666171116

$data[11,0] = 1100.0
$data[11,1] = "struggle"
$data[11,2] = -1.787458419799805
$data[11,3] = 1.655651211738586
$data[11,4] = -5.035046577453613
$data[11,5] = 0.5366638071683012
$data[11,6] = -2.253858975001746
$data[11,7] = 0.0269684557403839

$data[12,0] = 1200.0
$data[12,1] = "struggle"
$data[12,2] = -9.838252067565918
$data[12,3] = 3.984453201293945
$data[12,4] = -6.098217010498047
$data[12,5] = -0.129132547548839
$data[12,6] = -0.7719840942596894
$data[12,7] = 0.03628414990950613

$data[13,0] = 1300.0
$data[13,1] = "struggle"
$data[13,2] = 6.411758422851562
$data[13,3] = 1.583425164222717
$data[13,4] = 7.352428436279297
$data[13,5] = 0.3077981770038601
$data[13,6] = 1.624053824921043
$data[13,7] = 0.1965552446793536

$data[14,0] = 1400.0
$data[14,1] = "struggle"
$data[14,2] = -2.261712551116944
$data[14,3] = 0.8220813274383545
$data[14,4] = 1.315514087677002
$data[14,5] = 1.183255352536026
$data[14,6] = 3.720431172117896
$data[14,7] = 0.4516974523359422

$data[15,0] = 1500.0
$data[15,1] = "struggle"
$data[15,2] = -1.646389007568359
$data[15,3] = 0.2190679311752319
$data[15,4] = 0.9841623306274414
$data[15,5] = 0.7590655258723678
$data[15,6] = 3.928779942648755
$data[15,7] = 1.171988606452941

$data[16,0] = 1600.0
$data[16,1] = "struggle"
$data[16,2] = -1.105591297149658
$data[16,3] = 1.020219326019287
$data[16,4] = 3.201179504394531
$data[16,5] = 0.2291679642334281
$data[16,6] = 2.087977978647975
$data[16,7] = 1.354070066189282

$data[17,0] = 1700.0
$data[17,1] = "struggle"
$data[17,2] = -4.33466100692749
$data[17,3] = -0.8289146423339844
$data[17,4] = 6.12528133392334
$data[17,5] = 0.03184602683296017
$data[17,6] = 0.9907392433711465
$data[17,7] = 0.4831009315592913

$data[18,0] = 1800.0
$data[18,1] = "struggle"
$data[18,2] = -1.558335304260254
$data[18,3] = -0.159212052822113
$data[18,4] = 1.605715155601502
$data[18,5] = -0.04444044737183325
$data[18,6] = 0.1521864691559109
$data[18,7] = -0.2139835976520377

$data[19,0] = 1900.0
$data[19,1] = "struggle"
$data[19,2] = -0.9647946357727052
$data[19,3] = 1.00678539276123
$data[19,4] = -4.680802822113037
$data[19,5] = -0.4428928944529248
$data[19,6] = -2.088161782342548
$data[19,7] = -1.22169929499528

$data[20,0] = 2000.0
$data[20,1] = "struggle"
$data[20,2] = -3.810809135437012
$data[20,3] = 1.403007388114929
$data[20,4] = 0.0495486259460449
$data[20,5] = -1.581159264457474
$data[20,6] = -3.781517471585969
$data[20,7] = -2.15734222470499

$data[21,0] = 2100.0
$data[21,1] = "struggle"
$data[21,2] = -1.585423946380615
$data[21,3] = 2.060841083526612
$data[21,4] = -2.507726192474365
$data[21,5] = -0.5357818153439736
$data[21,6] = -0.6552340047700065
$data[21,7] = 0.5789350879435637

$data[22,0] = 2200.0
$data[22,1] = "struggle"
$data[22,2] = -5.486822128295898
$data[22,3] = 2.457437515258789
$data[22,4] = -1.076503276824951
$data[22,5] = -0.4312272305999488
$data[22,6] = -0.03695735122476339
$data[22,7] = -0.2086323031357354

$data[23,0] = 2300.0
$data[23,1] = "struggle"
$data[23,2] = 3.813155174255371
$data[23,3] = -5.157403945922852
$data[23,4] = 7.194998264312744
$data[23,5] = 0.1521366113910867
$data[23,6] = 0.3846518628451288
$data[23,7] = -0.3579327458021597

$data[24,0] = 2400.0
$data[24,1] = "struggle"
$data[24,2] = -3.507768154144287
$data[24,3] = 2.501498937606812
$data[24,4] = 0.7795240879058838
$data[24,5] = 0.3417635331956719
$data[24,6] = 0.5353018106246487
$data[24,7] = 0.2242374224018085

$data[25,0] = 2500.0
$data[25,1] = "struggle"
$data[25,2] = 0.2215757369995117
$data[25,3] = -0.4009582996368408
$data[25,4] = 2.163901329040528
$data[25,5] = -0.1421539567563001
$data[25,6] = 0.2858568746216445
$data[25,7] = -0.08254160519157175

$data[26,0] = 2600.0
$data[26,1] = "struggle"
$data[26,2] = 0.1625576019287109
$data[26,3] = 1.34720504283905
$data[26,4] = -0.6319388151168823
$data[26,5] = -0.06803667803808118
$data[26,6] = 0.1790894811250734
$data[26,7] = 0.08434615633925557

$data[27,0] = 2700.0
$data[27,1] = "struggle"
$data[27,2] = 0.044438362121582
$data[27,3] = -0.1398162841796875
$data[27,4] = -0.8414495587348938
$data[27,5] = -0.08848196070413178
$data[27,6] = -0.1032362286837731
$data[27,7] = 0.232543302129726

$data[28,0] = 2800.0
$data[28,1] = "struggle"
$data[28,2] = -0.1983919143676757
$data[28,3] = -0.413076639175415
$data[28,4] = 0.2017757892608642
$data[28,5] = 0.06249837318853448
$data[28,6] = 0.04699299066346516
$data[28,7] = 0.2259266389419836

$data[29,0] = 2900.0
$data[29,1] = "struggle"
$data[29,2] = 0.6235456466674805
$data[29,3] = 1.087465167045593
$data[29,4] = 0.6343502402305603
$data[29,5] = -0.0310014113783836
$data[29,6] = 0.0229074470698833
$data[29,7] = -0.0424551330506801

$ws.Range("A2:H31").Value = $data

# dimension / used range will auto-extend to A1:H31 once the new rows are written
